$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Notes" sheet
$ws2 = $wb.Worksheets.Item(2)   # "Data" sheet

# --- Notes sheet: add "Winter" to the Semesters list -----------------
# Insert a new blank row at 17 (between "Summer" at 16 and the blank
# separator row that was at 18) which pushes everything below it down
# by one row, then fill it with the new "Winter" entry.
[void]$ws1.Rows("17:17").Insert()
$ws1.Range("A17").Value = "Winter"

# Leave the Notes sheet's selection where Excel would land after the
# insert/edit.
[void]$ws1.Range("A22").Select()

# --- Data sheet: add data validation for nrows / ABET columns --------
# Whole-number validation (calendar year) on column E.
[void]$ws2.Range("E1:E1048576").Validation.Add(1, 1, 1, 1900, 2100)

# List validation on column C (Program Type) driven by the Notes sheet
# Program Types list.
[void]$ws2.Range("C1:C1048576").Validation.Add(3, 1, 3, "=Notes!A7:A11")
$ws2.Range("C1:C1048576").Validation.Formula1 = "=Notes!`$A`$7:`$A`$11"

# List validation on column D (Term) driven by the Notes sheet
# Semesters list (now Fall/Spring/Summer/Winter after the insert above).
[void]$ws2.Range("D1:D1048576").Validation.Add(3, 1, 3, "=Notes!A14:A17")
$ws2.Range("D1:D1048576").Validation.Formula1 = "=Notes!`$A`$14:`$A`$17"

# --- Switch the active tab from Notes to Data -------------------------
[void]$ws2.Activate()
[void]$ws2.Range("A2").Select()
